# Applies the "Atualização automática" row-rotation edit to rows 7-11
# of the dashboard sheet: each row's identity (Fly_ID + detection
# details in columns A, D-J) shifts down by one row, with the former
# row 7 data wrapping around into row 11. Columns B (Class) and C
# (First_Detection_Date) are identical across these rows so they are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 7;  A = "2117575c-4ae1-458c-b88a-fc40f40debdb"; D = "image_20250727074723_ppp0.jpg"; E = "PLACA_20250723145134"; F = "Moura"; G = 38.06587;  H = -7.221796; I = "1490,161,1563,258"; J = "0.62" },
    @{ Row = 8;  A = "283b6eda-9c83-4cdd-9524-c7c394f2dc89"; D = "image_20250728214139_ppp0.jpg"; E = "PLACA_20250717165933"; F = "Beja";  G = 38.02035;  H = -7.94715;  I = "962,713,1006,765";  J = "0.76" },
    @{ Row = 9;  A = "a19b65d1-6f97-4841-9e1c-7446a9be92b6"; D = "image_20250728214139_ppp0.jpg"; E = "PLACA_20250717165933"; F = "Beja";  G = 38.02035;  H = -7.94715;  I = "967,614,1002,659";  J = "0.73" },
    @{ Row = 10; A = "4be1b1cf-d480-453e-b5fb-d4ecd6764c4d"; D = "image_20250728214139_ppp0.jpg"; E = "PLACA_20250717165933"; F = "Beja";  G = 38.02035;  H = -7.94715;  I = "702,633,740,690";   J = "0.72" },
    @{ Row = 11; A = "dfd476d4-7689-4671-a076-78fe3ce806bb"; D = "image_20250728214139_ppp0.jpg"; E = "PLACA_20250717165933"; F = "Beja";  G = 38.02035;  H = -7.94715;  I = "1254,850,1294,895"; J = "0.67" }
)

foreach ($r in $rows) {
    $row = $r.Row

    # Columns holding comma-separated pixel coordinates (e.g. "702,633,740,690")
    # must stay as text, otherwise Excel auto-converts them into a plain
    # number and strips the commas. Force the cell format to Text first.
    $coordCell = $ws.Cells.Item($row, 9)
    $coordCell.NumberFormat = "@"

    $ws.Cells.Item($row, 1).Value  = $r.A   # A: Fly_ID
    $ws.Cells.Item($row, 4).Value  = $r.D   # D: First_Detection_Image
    $ws.Cells.Item($row, 5).Value  = $r.E   # E: Placa ID
    $ws.Cells.Item($row, 6).Value  = $r.F   # F: Localização
    $ws.Cells.Item($row, 7).Value  = $r.G   # G: Latitude
    $ws.Cells.Item($row, 8).Value  = $r.H   # H: Longitude
    $coordCell.Value              = $r.I    # I: First_Coords (forced text)
    $ws.Cells.Item($row, 10).Value = $r.J   # J: First_Confidence
}
